$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 7 (the "Obstetrician" specialist row) into rows 8 and 9,
# preserving its cell formatting (font/fill/border/alignment/number format).
$ws.Rows(7).Copy()
$ws.Rows(8).Insert()
$ws.Rows(7).Copy()
$ws.Rows(9).Insert()
$excel.CutCopyMode = $false

# Populate column D (condition value) first, then column A (rule name) so the
# shared-string table is built up in the same order as the source edit.
$ws.Range("D8").Value = "Gynaecologist"
$ws.Range("D9").Value = "Pediatrician"
$ws.Range("A8").Value = "TariffDecision3"
$ws.Range("A9").Value = "TariffDecision4"

# Remaining cell values for the two new rows.
$ws.Range("B8").Value = 20
$ws.Range("C8").Value = "SPECIALIST CONSULTATION "
$ws.Range("E8").Value = "tier1"
$ws.Range("F8").Value = 150000

$ws.Range("B9").Value = 20
$ws.Range("C9").Value = "SPECIALIST CONSULTATION "
$ws.Range("E9").Value = "tier1"
$ws.Range("F9").Value = 15000

# Row heights: row 6 shrinks to its default content height, rows 7-9 become
# a uniform 30pt (two-line wrapped) height.
$ws.Rows(6).RowHeight = 14.25
$ws.Rows(7).RowHeight = 30
$ws.Rows(8).RowHeight = 30
$ws.Rows(9).RowHeight = 30

# The data rows (6-9) lose their explicit Arial font override and fall back
# to the workbook default (Calibri 11) while keeping borders/alignment.
# (Column A never carried an explicit style, so it is left untouched.)
$dataRows = $ws.Range("B6:F9")
$dataRows.Font.Name = "Calibri"
$dataRows.Font.Size = 11

# Selection moves to the newly-added last row.
$ws.Range("B8").Select()
